# This script applies the "ws-01" RDBMS saveResults / sound-category commit to the
# hidden '#system' sheet that backs the workbook's defined-name command catalog.
#
# Summary of the edit:
#  - A brand new "sound" command category is inserted as a new column (Q) in the
#    '#system' sheet, between "redis" (col P) and the former "ssh" column, which
#    (along with every category to its right: ssh, step, web, webalert, webcookie,
#    ws, xml) shifts one column to the right.
#  - "sound" is also inserted (alphabetically) into the list of categories kept in
#    column A ("target"), and it gets two commands of its own: lazer(repeats) and
#    warp(repeats).
#  - A new "clear(vars)" command is inserted (alphabetically) into the "base"
#    category (column C).
#  - A new "saveRowCount(var)" command is inserted (alphabetically) into the
#    "desktop" category (column E).
#  - A new "saveResults(db,sqls,outputDir)" command is appended to the end of the
#    "rdbms" category (column O).
#  - Two existing command signatures are updated in place: mail's send(...) gains a
#    "profile" parameter, and number's assertBetween(...) renames lower/upper to
#    min/max.
#  - All affected defined names are updated to the new ranges, and a new "sound"
#    defined name is added.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

function Set-DefinedName($name, $refersTo) {
    $wb.Names.Item($name).RefersTo = $refersTo
}

# --- 1. Insert a brand-new column at Q ("sound"), shifting old Q..W (ssh..xml) one
#        column to the right, into R..X. This is a genuine whole-column insert, so
#        Excel's native Insert/ShiftToRight behavior is exactly what's needed here. ---
$ws.Columns("Q").Insert(-4161)  # xlShiftToRight

# --- 2. "target" (column A): insert a new row at A17 for "sound", pushing the
#        previously-17th..23rd entries (ssh..xml) down to 18..24. Only column A
#        should move, so this is done manually cell-by-cell rather than via Insert
#        (which shifts an entire worksheet row here). ---
for ($i = 23; $i -ge 17; $i--) {
    $ws.Cells.Item($i + 1, 1).Value = $ws.Cells.Item($i, 1).Value2
}
$ws.Cells.Item(17, 1).Value = "sound"

# --- 3. "base" (column C): insert a new row at C16 for "clear(vars)", pushing the
#        previously-16th..31st entries down to 17..32. ---
for ($i = 31; $i -ge 16; $i--) {
    $ws.Cells.Item($i + 1, 3).Value = $ws.Cells.Item($i, 3).Value2
}
$ws.Cells.Item(16, 3).Value = "clear(vars)"

# --- 4. "desktop" (column E): insert a new row at E69 for "saveRowCount(var)",
#        pushing the previously-69th..90th entries down to 70..91. ---
for ($i = 90; $i -ge 69; $i--) {
    $ws.Cells.Item($i + 1, 5).Value = $ws.Cells.Item($i, 5).Value2
}
$ws.Cells.Item(69, 5).Value = "saveRowCount(var)"

# --- 5. "rdbms" (column O): append the new saveResults(...) command as the 7th
#        entry (row 7); no shifting needed since it sorts to the end. ---
$ws.Cells.Item(7, 15).Value = "saveResults(db,sqls,outputDir)"

# --- 6. Update two existing command signatures in place. ---
$ws.Range("L2").Value = "send(profile,to,subject,body)"
$ws.Range("M2").Value = "assertBetween(num,min,max)"

# --- 7. Populate the new "sound" column (Q) header + its two commands. ---
$ws.Range("Q1").Value = "sound"
$ws.Range("Q2").Value = "lazer(repeats)"
$ws.Range("Q3").Value = "warp(repeats)"

# --- 8. Update every defined name whose range moved, and add the new one. ---
Set-DefinedName "base" "='#system'!`$C`$2:`$C`$32"
Set-DefinedName "desktop" "='#system'!`$E`$2:`$E`$91"
Set-DefinedName "rdbms" "='#system'!`$O`$2:`$O`$7"
Set-DefinedName "ssh" "='#system'!`$R`$2:`$R`$9"
Set-DefinedName "target" "='#system'!`$A`$2:`$A`$24"
Set-DefinedName "web" "='#system'!`$T`$2:`$T`$108"
Set-DefinedName "webalert" "='#system'!`$U`$2:`$U`$6"
Set-DefinedName "webcookie" "='#system'!`$V`$2:`$V`$8"
Set-DefinedName "ws" "='#system'!`$W`$2:`$W`$16"
Set-DefinedName "xml" "='#system'!`$X`$2:`$X`$11"
Set-DefinedName "step" "='#system'!`$S`$2:`$S`$4"
$wb.Names.Add("sound", "='#system'!`$Q`$2:`$Q`$3")
